$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: the in-progress entry now has an end time recorded ---
$ws.Range("D31").Value = 0.6875
$ws.Range("D31").NumberFormat = $ws.Range("D30").NumberFormat

# --- Rows 29-31: extend the "end minus start" column down (shared formula) ---
$ws.Range("E29:E31").Formula = "=D29-C29"
$ws.Range("E29:E31").NumberFormat = $ws.Range("E28").NumberFormat

# --- New day started: row 33 holds the date, row 32 stays blank ---
$ws.Range("B33").Value = 42929
$ws.Range("B30").Copy()
$ws.Range("B33").PasteSpecial(-4122)

# --- Column D grew a bit wider to fit the new values ---
$ws.Columns("D").AutoFit()

# --- Put the selection where the user left off editing ---
$ws.Range("C33").Select()
